$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "a"
$ws.Range("A3").Value = "b1"
$ws.Range("A4").Value = "b2"
$ws.Range("A5").Value = "c1"
$ws.Range("A6").Value = "c2"
$ws.Range("A7").Value = "c3"
$ws.Range("A8").Value = "c4"
